# Append 1 row (row 5) of data to Sheet1 - matches the commit
# "Append 1 rows at 2025-05-01T16:36:50.762Z".
#
# Column A is blank and column C holds a number-looking value ("222");
# both need to be stored as literal TEXT (like the other rows in this
# sheet), not as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

# Leading apostrophe forces Excel to store the (empty) entry as literal text,
# matching the existing blank-but-text A2:A4 cells in this column.
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 2).Value = "يامن "

# Force text storage so "222" isn't reinterpreted as a number.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "222"

$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"
$ws.Cells.Item($row, 6).Value = "C1"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٧:٣٦:٥٠ م"
